$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.092.17"
$ws.Range("E2").Value = "  -0.83%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.650.28"
$ws.Range("E3").Value = "  -0.90%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.32"
$ws.Range("E5").Value = "  -0.85%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5202"
$ws.Range("E6").Value = "  -2.74%  "

$ws.Range("E7").Value = "  -0.38%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2614"
$ws.Range("E8").Value = "  -1.78%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06281"
$ws.Range("E9").Value = "  -2.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.50"
$ws.Range("E10").Value = "  -0.53%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07800"
$ws.Range("E11").Value = "  -0.23%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.703.78"
$ws.Range("E12").Value = "  +2.52%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.459"
$ws.Range("E13").Value = "  -2.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.877.63"
$ws.Range("E14").Value = "  -0.77%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5540"
$ws.Range("E15").Value = "  +0.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅7987"
$ws.Range("E16").Value = "  -2.62%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.68"
$ws.Range("E17").Value = "  -1.72%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.079.94"
$ws.Range("E18").Value = "  -0.91%  "

$ws.Range("E19").Value = "  -0.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.627"
$ws.Range("E20").Value = "  -1.45%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.92"
$ws.Range("E21").Value = "  -0.08%  "

$ws.Range("E22").Value = "  -1.88%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.939"
$ws.Range("E23").Value = "  -1.75%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.006"
$ws.Range("E24").Value = "  -0.49%  "

$ws.Range("E25").Value = "  +0.29%  "

$ws.Range("E26").Value = "  -2.42%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.180"
$ws.Range("E27").Value = "  -0.28%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.87"
$ws.Range("E28").Value = "  -1.49%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.476"
$ws.Range("E29").Value = "  -0.56%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05601"
$ws.Range("E30").Value = "  -3.96%  "

$ws.Range("E31").Value = "  -1.38%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.478"
$ws.Range("E32").Value = "  -3.79%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.363"
$ws.Range("E33").Value = "  +2.25%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.596"
$ws.Range("E34").Value = "  -1.06%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.801"
$ws.Range("E35").Value = "  -0.89%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9474"
$ws.Range("E36").Value = "  -2.05%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.405"
$ws.Range("E37").Value = "  -0.60%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5641"
$ws.Range("E38").Value = "  -2.91%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.973"
$ws.Range("E39").Value = "  +1.73%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01580"
$ws.Range("E40").Value = "  -1.57%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.058.81"
$ws.Range("E41").Value = "  +0.52%  "

$ws.Range("E42").Value = "  -0.49%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8381"
$ws.Range("E43").Value = "  -3.57%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.45"
$ws.Range("E44").Value = "  -2.31%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.788.14"
$ws.Range("E45").Value = "  -0.84%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.04"
$ws.Range("E46").Value = "  -1.62%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₈106"
$ws.Range("E47").Value = "  -0.46%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05380"
$ws.Range("E48").Value = "  +4.14%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.006"
$ws.Range("E49").Value = "  -0.75%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4333"
$ws.Range("E50").Value = "  -1.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.943"
$ws.Range("E51").Value = "  -1.17%  "
